$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: copy formatting (bold/border/center) from an existing
# header cell (AC1) so the new headers match the look of the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-49) gets the same team record: 79 wins, 83 losses, 0 ties.
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 79
    $ws.Cells.Item($r, 31).Value = 83
    $ws.Cells.Item($r, 32).Value = 0
}
